$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "write group": put the new group guid "123" into A3 (Group guid column).
# Entering it as a literal would auto-convert to a number, so build it as a
# text-formula result first, then paste-special (values only) over itself.
# That leaves a plain shared-string text cell without touching the cell's
# number format / style (stays on the default style, like in the target).
$ws.Range("A3").Formula = "=""123"""
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# "write material": put the new group name into D3 (group name column).
$ws.Range("D3").Value = "Майки / Футболки / Рубашки / Комбинезоны"
